# Updated symbol list on Wed Dec 21 11:29:06 UTC 2022 with GitHub Actions
# Refresh cryptocurrency price/volume snapshot values on Sheet1.
# Price values in column D are numeric-looking text, so they are entered
# with a leading apostrophe to force Excel to keep them as text (preserving
# exact formatting / trailing zeros), then the quote-prefix style is reset
# back to Normal so no extra cell styling is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'248.91"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").Value = "'22.70"
$ws.Range("D3").Style = "Normal"

$ws.Range("D4").Value = "'5.377"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").Value = "'0.05686"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").Value = "'3.404"
$ws.Range("D6").Style = "Normal"

$ws.Range("D8").Value = "'0.8059"
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").Value = "'0.9152"
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").Value = "'0.1406"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").Value = "'0.07442"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").Value = "'0.03106"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").Value = "'0.03052"
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").Value = "'0.09383"
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").Value = "'3.874"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").Value = "'0.001582"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").Value = "'0.04799"
$ws.Range("D17").Style = "Normal"

$ws.Range("D20").Value = "'0.006494"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").Value = "'0.004996"
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").Value = "'0.0009988"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").Value = "'0.0001500"
$ws.Range("D23").Style = "Normal"

$ws.Range("D40").Value = "'0.03997"
$ws.Range("D40").Style = "Normal"

$ws.Range("B41").Value = "KickToken"

$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"

$ws.Range("D41").Value = "'0.006841"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"

$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"

$ws.Range("D42").Value = "'0.1069"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"

$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"

$ws.Range("D43").Value = "'0.002756"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").Value = "'0.007964"
$ws.Range("D44").Style = "Normal"

$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"

$ws.Range("D48").Value = "'0.2086"
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D49").Style = "Normal"
